# Weekly update: insert a new price record at row 33 for
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Alcachofa".
# All existing records from row 33 downward shift down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 33:83 down to 34:84, inserting a fresh blank row 33.
$ws.Rows.Item(33).Insert(-4121)

# Populate the brand-new row 33 with this week's record.
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 45079
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112013
$ws.Range("G33").Value = "Alcachofa"
$ws.Range("H33").Value = "Española"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 30
$ws.Range("K33").Value = 16000
$ws.Range("L33").Value = 16000
$ws.Range("M33").Value = 16000
$ws.Range("N33").Value = "`$/caja 30 unidades"
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 533
$ws.Range("Q33").Value = 30
$ws.Range("R33").Value = "Hortaliza"
